$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-10 (columns E..T)
$data = @{
    2  = @{ E=3; F=1; G=4.055887666666667; H=12.167663; I=0.4763357569530485; J=0.4763357569530485; K=3; L=1; M=0.380615; N=1.141845; O=0.04044104717146424; P=0.04044104717146424; Q=1.543731684248333; R=13.893585158235; S=0.01926351681639336; T=0.01926351681639336 }
    3  = @{ E=3; F=1; G=4.055887666666667; H=12.167663; I=0.4763357569530485; J=0.4763357569530485; K=3; L=1; M=4.865925666666667; N=14.597777; O=0.5170135948885495; P=0.5170135948885495; Q=19.73564789835011; R=177.620831085151; S=0.246272062076254; T=0.246272062076254 }
    4  = @{ E=3; F=1; G=4.055887666666667; H=12.167663; I=0.4763357569530485; J=0.4763357569530485; K=3; L=1; M=4.165060333333333; N=12.495181; O=0.4425453579399863; P=0.4425453579399863; Q=16.89301683688922; R=152.037151532003; S=0.2108001780604012; T=0.2108001780604012 }
    5  = @{ E=3; F=1; G=3.886460333333333; H=11.659381; I=0.4564376967244237; J=0.4564376967244237; K=3; L=1; M=0.380615; N=1.141845; O=0.04044104717146424; P=0.04044104717146424; Q=1.479245099771667; R=13.313205897945; S=0.01845881842406691; T=0.01845881842406691 }
    6  = @{ E=3; F=1; G=3.886460333333333; H=11.659381; I=0.4564376967244237; J=0.4564376967244237; K=3; L=1; M=4.865925666666667; N=14.597777; O=0.5170135948885495; P=0.5170135948885495; Q=18.91122708844856; R=170.201043796037; S=0.2359844944261438; T=0.2359844944261438 }
    7  = @{ E=3; F=1; G=3.886460333333333; H=11.659381; I=0.4564376967244237; J=0.4564376967244237; K=3; L=1; M=4.165060333333333; N=12.495181; O=0.4425453579399863; P=0.4425453579399863; Q=16.18734177144011; R=145.686075942961; S=0.201994383874213; T=0.201994383874213 }
    8  = @{ E=3; F=1; G=0.5724183333333334; H=1.717255; I=0.06722654632252778; J=0.06722654632252777; K=3; L=1; M=0.380615; N=1.141845; O=0.04044104717146424; P=0.04044104717146424; Q=0.2178710039416667; R=1.960839035475; S=0.002718711931003972; T=0.002718711931003971 }
    9  = @{ E=3; F=1; G=0.5724183333333334; H=1.717255; I=0.06722654632252778; J=0.06722654632252777; K=3; L=1; M=4.865925666666667; N=14.597777; O=0.5170135948885495; P=0.5170135948885495; Q=2.785345060237222; R=25.068105542135; S=0.03475703838615168; T=0.03475703838615168 }
    10 = @{ E=3; F=1; G=0.5724183333333334; H=1.717255; I=0.06722654632252778; J=0.06722654632252777; K=3; L=1; M=4.165060333333333; N=12.495181; O=0.4425453579399863; P=0.4425453579399863; Q=2.384156894239444; R=21.457412048155; S=0.02975079600537213; T=0.02975079600537212 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
